# Enable geofence export for PDF and Excel logbook template.
# Inserts two new columns (N, O) before the existing "Notiz" column (which
# shifts from N to P), adds "Start Geofence" / "Ende Geofences" headers and
# their corresponding ${entry.startGeofences} / ${entry.endGeofences}
# template placeholders, updates the jxls comments' lastCell references from
# N13 to P13, and updates the saved cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at N:O - this shifts the existing "Notiz" column
# (N) and its formatting out to column P, preserving styles/column widths.
$ws.Columns("N:O").Insert()

# Populate the new cells in the same order the original workbook's shared
# strings were appended (N12, N13, O12, O13) so the shared-string table
# indices line up with the authored change.
$ws.Range("N12").Value = "Start Geofence"
$ws.Range("N13").Value = '${entry.startGeofences}'
$ws.Range("O12").Value = "Ende Geofences"
$ws.Range("O13").Value = '${entry.endGeofences}'

# Update the jxls directives stored in cell comments so the export area /
# loops now extend through column P instead of N.
foreach ($cmt in $ws.Comments) {
    $text = $cmt.Text()
    $updated = $text.Replace('lastCell="N13"', 'lastCell="P13"')
    [void]$cmt.Text($updated)
}

# Restore/update the active cell selection recorded with the sheet view.
[void]$ws.Range("F31").Select()
